$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.715.67"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.233.95"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.01"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.06"
$ws.Range("E6").Value = "  +4.66%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.233.40"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.69"
$ws.Range("E11").Value = "  -6.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.511"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  +4.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.16"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "3.766.84"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "66.759.22"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.49"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "3.239.98"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "512.93"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.35"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.09"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.91"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.00"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.01"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.26"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.41"
$ws.Range("E29").Value = "  +6.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.03"
$ws.Range("E30").Value = "  +6.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.02"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.30"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.57"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "516.95"
$ws.Range("E36").Value = "  +7.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0960"
$ws.Range("E37").Value = "  +7.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.15"
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("D39").Value = "0.0₃0778"
$ws.Range("E39").Value = "  +22.50%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.10"
$ws.Range("E40").Value = "  +9.71%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0422"
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("E42").Value = "  +6.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.82"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.302"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("D46").Value = "2.869.31"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.70"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +6.03%  "
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.63"
$ws.Range("E51").Value = "  +1.90%  "
